$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (header "K", previously holding raw Strike# totals) with
# recomputed K values per row, per commit: "regen save_data to use K instead
# of Strike#, regen std/mean, calc and write s_vals"
$gValues = @{
    2 = 3
    3 = 3
    4 = 2
    5 = 1
    6 = 0
    7 = 2
    8 = 0
    9 = 3
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 2
    20 = 2
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 0
    29 = 0
    30 = 0
    31 = 5
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 3
    41 = 2
    42 = 0
    43 = 3
    44 = 1
    45 = 8
    46 = 4
    47 = 2
    48 = 5
    49 = 6
    50 = 8
    51 = 4
    52 = 5
    53 = 4
    54 = 3
    55 = 6
    56 = 7
    57 = 4
    58 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

